$d = $word.ActiveDocument

function New-WordXmlPackage($bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Build the 4 new paragraphs that need to be inserted before the first
# existing paragraph ("How much overhead is incurred by a function call?").
$newParasXml = '<w:body>' +
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Page Numbers of Some Important Figures</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Page 289</w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  GDB commands</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>' +
  '</w:body>'

# The trailing (last) paragraph currently holds the "_GoBack" bookmark.
# In the edited document that bookmark is relocated into one of the new
# paragraphs being inserted at the top, so remove it here first, while it
# is still the only "_GoBack" bookmark in the document (avoids any
# ambiguity once the new bookmark is inserted below).
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

$firstPara = $d.Paragraphs.Item(1)
$insertRange = $firstPara.Range
$insertRange.Collapse(1)   # wdCollapseStart
[void]$insertRange.InsertXML((New-WordXmlPackage $newParasXml))
